# Activity_Utilization.xlsx - "more test cases completed"
# Adds new negative-test-case rows (row 3) to the GET_last_login,
# GET_equipment_session and GET_equipment_summary sheets, fixes the
# description of the existing GET_equipment_summary test case, and
# tweaks a couple of column widths / selections left behind by the
# author while working in Excel.

$wb = $excel.ActiveWorkbook

$wsEquipment = $wb.Worksheets.Item("GET_equipment")
$wsLastLogin = $wb.Worksheets.Item("GET_last_login")
$wsSession   = $wb.Worksheets.Item("GET_equipment_session")
$wsSummary   = $wb.Worksheets.Item("GET_equipment_summary")

function Set-TextCellLikeNeighbor($DstCell, $SrcCell, $Value) {
    $DstCell.Font.Name = $SrcCell.Font.Name
    $DstCell.Font.Size = $SrcCell.Font.Size
    $DstCell.NumberFormat = $SrcCell.NumberFormat
    $DstCell.Value = $Value
}

# ---------------------------------------------------------------------
# GET_equipment (sheet1): only cosmetic leftovers - wider Description
# column and a different cell selected when the author saved the file.
# ---------------------------------------------------------------------
$wsEquipment.Columns.Item(3).ColumnWidth = 32.166666666666664
$wsEquipment.Activate()
$wsEquipment.Range("D16").Select()

# ---------------------------------------------------------------------
# GET_last_login (sheet2): new negative test case - invalid equipmentId
# GET_equipment_session (sheet3): new negative test case
#
# The two sheets were edited together (Description/schemaValidationFile/
# Uri filled in for sheet2, then the same columns for sheet3, then the
# ExpectedStatusCode column was filled in for both) which is why the
# shared-string table interleaves them the way it does.
# ---------------------------------------------------------------------
$wsLastLogin.Range("A3").Value = 2
$wsLastLogin.Range("B3").Value = "Y"
Set-TextCellLikeNeighbor $wsLastLogin.Range("C3") $wsLastLogin.Range("C2") "Get Last Logins by invalid Equipment ID"
Set-TextCellLikeNeighbor $wsLastLogin.Range("H3") $wsLastLogin.Range("H2") "400error.json"
Set-TextCellLikeNeighbor $wsLastLogin.Range("F3") $wsLastLogin.Range("F2") "/activity/v1/equipment/bfc105b5-7ae6-432d-ae13-db2a3096a0c9/last-logins?pageSize=50&startTimestamp=2021-5-1&endTimestamp=2022-5-1"
Set-TextCellLikeNeighbor $wsLastLogin.Range("E3") $wsLastLogin.Range("E2") "GET"

$wsSession.Range("A3").Value = 2
$wsSession.Range("B3").Value = "Y"
Set-TextCellLikeNeighbor $wsSession.Range("C3") $wsSession.Range("C2") "Get equipment session with invalid equipmentId"
Set-TextCellLikeNeighbor $wsSession.Range("F3") $wsSession.Range("F2") "/activity/v1/equipment/bfc105b5-7ae6-432d-ae13-db2a3096a0c9/sessions?pageSize=50&startTimestamp=2021-5-1&endTimestamp=2022-5-1"
Set-TextCellLikeNeighbor $wsSession.Range("H3") $wsSession.Range("H2") "400error.json"
Set-TextCellLikeNeighbor $wsSession.Range("E3") $wsSession.Range("E2") "GET"

Set-TextCellLikeNeighbor $wsLastLogin.Range("G3") $wsLastLogin.Range("G2") "400"
Set-TextCellLikeNeighbor $wsSession.Range("G3") $wsSession.Range("G2") "400"

$wsLastLogin.Activate()
$wsLastLogin.Range("G3").Select()

$wsSession.Activate()
$wsSession.Range("F3").Select()

# ---------------------------------------------------------------------
# GET_equipment_summary (sheet4): fix existing description, add new
# negative test case.
# ---------------------------------------------------------------------
$wsSummary.Columns.Item(3).ColumnWidth = 39.833333333333336

$wsSummary.Range("A3").Value = 2
$wsSummary.Range("B3").Value = "Y"
Set-TextCellLikeNeighbor $wsSummary.Range("F3") $wsSummary.Range("F2") "/activity/v1/equipment/bfc105b5-7ae6-432d-ae13-db2a3096a0c9/summary?pageSize=50&startTimestamp=2021-5-1&endTimestamp=2022-5-1"
$wsSummary.Range("E3").Value = "GET"
$wsSummary.Range("G3").Value = 400
$wsSummary.Range("H3").Value = "400error.json"

Set-TextCellLikeNeighbor $wsSummary.Range("C2") $wsSummary.Range("C2") "Get Equipment Summary"
Set-TextCellLikeNeighbor $wsSummary.Range("C3") $wsSummary.Range("C2") "Get Equipment Summary with invalid equipmentId"

$wsSummary.Activate()
$wsSummary.Range("C3").Select()
